$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.540407678265278
$ws.Range("B2").Value = -4.858687610135589

$ws.Range("A3").Value = -0.519008676897071
$ws.Range("B3").Value = 0.9505066755281615

$ws.Range("A4").Value = 0.9694740348354314
$ws.Range("B4").Value = -3.393773433574415

$ws.Range("A5").Value = 0.7134754889739495
$ws.Range("B5").Value = 0.8263685032558445

$ws.Range("A6").Value = -0.8211782141201494
$ws.Range("B6").Value = -1.800952556809995
